$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new location mapping ("Law Library > Special Reserve") was added right
# after "Law Library Technical Services" (row 46) and before
# "Law Library Legal Aid Clinic" (previously row 47). Insert a new row so
# everything below shifts down by one, then populate it.
$ws.Rows.Item(47).Insert()

$ws.Cells.Item(47, 1).Value = "Law Library Special Reserve"
$ws.Cells.Item(47, 5).Value = "Law Library > Special Reserve"

# Match the author's final cursor/selection position.
$ws.Range("A48").Select()
